$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J ("Strato 127") header + scattered data points ---
$ws.Range("J1").Value = "Strato 127"
$ws.Range("J2").Value = 19884
$ws.Range("J3").Value = 56897
$ws.Range("J4").Value = 59294
$ws.Range("J10").Value = 49978
$ws.Range("J21").Value = 59728

# The "sqrt of num computations" note on row 24 shifts from column J to column K
# now that J is used for the new "Strato 127" series.
$ws.Range("J24").ClearContents() | Out-Null
$ws.Range("K24").Value = "sqrt of num computations"

# --- New data rows 26-31 ---
$ws.Range("A26").Value = 100
$ws.Range("B26").Value = 100000
$ws.Range("C26").Value = $false
$ws.Range("D26").Value = $false
$ws.Range("E26").Value = 36117
$ws.Range("H26").Value = 36424

$ws.Range("A27").Value = 150
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = $false
$ws.Range("D27").Value = $false
$ws.Range("E27").Value = 61840

$ws.Range("A28").Value = 150
$ws.Range("B28").Value = 4
$ws.Range("C28").Value = $false
$ws.Range("D28").Value = $false
$ws.Range("E28").Value = 142366

$ws.Range("A29").Value = 200
$ws.Range("B29").Value = 0
$ws.Range("C29").Value = $false
$ws.Range("D29").Value = $true
$ws.Range("E29").Value = "Out of heap space error"
$ws.Range("J29").Value = 1105204

$ws.Range("A30").Value = 200
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = $false
$ws.Range("D30").Value = $false
$ws.Range("E30").Value = "Out of heap space error"
$ws.Range("J30").Value = 379378

$ws.Range("A31").Value = 250
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = $false
$ws.Range("D31").Value = $false
$ws.Range("J31").Value = "Out of heap space"

# --- Column width tweaks (E got wider to fit the new 6-digit value on row 28;
#     J is a brand-new column sized to fit its content) ---
$ws.Columns("E").ColumnWidth = 6.25
$ws.Columns("J").ColumnWidth = 8.6

# --- Selection moves to the new first empty row beneath the data ---
$ws.Range("A32").Select() | Out-Null
